$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cmaes")

# Row 1: B1, D1, E1 change from "TVM" to "TVH" (C1 already "TVH")
$ws.Range("B1").Value = "TVH"
$ws.Range("D1").Value = "TVH"
$ws.Range("E1").Value = "TVH"

# Row 3: rename the "beats_simulation" knob to "scenario_ptr" and change the
# formula string used to load the scenario; also drop the duplicated values
# that used to live in C3:E3 (the "beats_simulation" concept collapses to a
# single B3 value now).
$ws.Range("A3").Value = "scenario_ptr"
$ws.Range("B3").Value = 'Utilities.char2char(''C:\Users\Felix\code\autoCalibrationProject\beats_scenarios\210E.xml'');'
$ws.Range("C3:E3").ClearContents()

# Rows 4-8: the C:E columns no longer carry duplicated per-run values, only
# column B remains populated.
$ws.Range("C4:E4").ClearContents()
$ws.Range("C5:E5").ClearContents()
$ws.Range("C6:E6").ClearContents()
$ws.Range("C7:E7").ClearContents()
$ws.Range("C8:E8").ClearContents()

# Update the active selection to C2 (was B23).
$ws.Activate()
$ws.Range("C2").Select()
